$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("K5").Value = 9
$ws.Range("L8").Value = 1.33
$ws.Range("M8").Value = 3.25
$ws.Range("G10").Value = 1.45
$ws.Range("H10").Value = 3.8
$ws.Range("I10").Value = 9
$ws.Range("J10").Value = 1.1
$ws.Range("K10").Value = 7
$ws.Range("P10").Value = 1.5
$ws.Range("Q10").Value = 2.5
$ws.Range("W10").Value = 9
$ws.Range("X10").Value = 15
$ws.Range("Z10").Value = 7
$ws.Range("AA10").Value = 8
$ws.Range("AE10").Value = 15
$ws.Range("G11").Value = 1.85
$ws.Range("H11").Value = 3.4
$ws.Range("I11").Value = 4.33
$ws.Range("J11").Value = 1.07
$ws.Range("K11").Value = 9
$ws.Range("U11").Value = 8
$ws.Range("W11").Value = 15
$ws.Range("Z11").Value = 7.5
$ws.Range("AA11").Value = 6.5
$ws.Range("AD11").Value = 1250
$ws.Range("AE11").Value = 10
$ws.Range("AF11").Value = 21
$ws.Range("AG11").Value = 15
$ws.Range("AJ11").Value = 41
$ws.Range("G12").Value = 1.5
$ws.Range("H12").Value = 3.8
$ws.Range("I12").Value = 6.1
$ws.Range("L12").Value = 1.36
$ws.Range("M12").Value = 2.67
$ws.Range("N12").Value = 2.05
$ws.Range("O12").Value = 1.62
$ws.Range("P12").Value = 1.45
$ws.Range("Q12").Value = 2.4
$ws.Range("R12").Value = 2.18
$ws.Range("S12").Value = 1.53
$ws.Range("T12").Value = 5.3
$ws.Range("U12").Value = 6
$ws.Range("V12").Value = 8.75
$ws.Range("W12").Value = 9.75
$ws.Range("X12").Value = 14.5
$ws.Range("Y12").Value = 40
$ws.Range("Z12").Value = 8.25
$ws.Range("AB12").Value = 24
$ws.Range("AG12").Value = 21
$ws.Range("AH12").Value = 150
$ws.Range("AI12").Value = 80
$ws.Range("AJ12").Value = 100
$ws.Range("K13").Value = 10
$ws.Range("G16").Value = 2
$ws.Range("H16").Value = 3.7
$ws.Range("J16").Value = 1.03
$ws.Range("K16").Value = 15
$ws.Range("L16").Value = 1.18
$ws.Range("M16").Value = 4.5
$ws.Range("N16").Value = 1.62
$ws.Range("O16").Value = 2.25
$ws.Range("Q16").Value = 3.4
$ws.Range("R16").Value = 1.57
$ws.Range("S16").Value = 2.25
$ws.Range("T16").Value = 10
$ws.Range("Y16").Value = 21
$ws.Range("Z16").Value = 15
$ws.Range("AA16").Value = 7
$ws.Range("AD16").Value = 126
$ws.Range("AE16").Value = 13
$ws.Range("AF16").Value = 21
$ws.Range("G17").Value = 1.38
$ws.Range("I17").Value = 8
$ws.Range("N17").Value = 1.65
$ws.Range("O17").Value = 2.2
$ws.Range("P17").Value = 1.3
$ws.Range("R17").Value = 1.95
$ws.Range("S17").Value = 1.8
$ws.Range("W17").Value = 9
$ws.Range("Y17").Value = 26
$ws.Range("Z17").Value = 13
$ws.Range("AE17").Value = 21
$ws.Range("G19").Value = 1.95
$ws.Range("H19").Value = 3.5
$ws.Range("I19").Value = 3.5
$ws.Range("J19").Value = 1.04
$ws.Range("K19").Value = 9
$ws.Range("N19").Value = 1.75
$ws.Range("O19").Value = 2.05
$ws.Range("U19").Value = 10
$ws.Range("W19").Value = 17
$ws.Range("AH19").Value = 41
$ws.Range("AI19").Value = 29
$ws.Range("G20").Value = 1.17
$ws.Range("H20").Value = 7.5
$ws.Range("L20").Value = 1.07
$ws.Range("M20").Value = 7.5
$ws.Range("N20").Value = 1.25
$ws.Range("O20").Value = 3.75
$ws.Range("P20").Value = 1.17
$ws.Range("Q20").Value = 4.5
$ws.Range("T20").Value = 13
$ws.Range("V20").Value = 11
$ws.Range("W20").Value = 8.5
$ws.Range("Z20").Value = 29
$ws.Range("AB20").Value = 23
$ws.Range("AD20").Value = 151
$ws.Range("AF20").Value = 51
